$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first data row (old row 2); this shifts the old row 3 up into
# row 2's place, matching the target layout exactly.
$ws.Rows(2).Delete()
